$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.05019399380407
$ws.Range("D2").Value = 1.05783290532045
$ws.Range("E2").Value = 1.057187736487439
$ws.Range("F2").Value = 1.067625180481724
$ws.Range("I2").Value = 1.049160186023271
$ws.Range("J2").Value = 1.055228439264246
$ws.Range("K2").Value = 1.060566794428925
$ws.Range("L2").Value = 1.05992339209416
$ws.Range("M2").Value = 1.070332541542478
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051099140355649
$ws.Range("D3").Value = 1.058565776587103
$ws.Range("E3").Value = 1.057986622290976
$ws.Range("F3").Value = 1.068490187668519
$ws.Range("I3").Value = 1.04942122303769
$ws.Range("J3").Value = 1.055783386719374
$ws.Range("K3").Value = 1.061113785461762
$ws.Range("L3").Value = 1.06053610267233
$ws.Range("M3").Value = 1.071013247691224
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051685383885528
$ws.Range("D4").Value = 1.059040430719538
$ws.Range("E4").Value = 1.058504402466647
$ws.Range("F4").Value = 1.06905081019511
$ws.Range("I4").Value = 1.049589069366782
$ws.Range("J4").Value = 1.056142366348279
$ws.Range("K4").Value = 1.06146749172106
$ws.Range("L4").Value = 1.06093275705777
$ws.Range("M4").Value = 1.07145396915327
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051931971740138
$ws.Range("D5").Value = 1.059240078431378
$ws.Range("E5").Value = 1.058722278545593
$ws.Range("F5").Value = 1.069286710867087
$ws.Range("I5").Value = 1.049659376964507
$ws.Range("J5").Value = 1.056293254248416
$ws.Range("K5").Value = 1.061616132702163
$ws.Range("L5").Value = 1.061099554678211
$ws.Range("M5").Value = 1.071639309055422
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.05197338257282
$ws.Range("D6").Value = 1.059273606160731
$ws.Range("E6").Value = 1.058758872654368
$ws.Range("F6").Value = 1.069326332183812
$ws.Range("I6").Value = 1.0496711669449
$ws.Range("J6").Value = 1.056318587365683
$ws.Range("K6").Value = 1.061641086791156
$ws.Range("L6").Value = 1.061127563282144
$ws.Range("M6").Value = 1.071670431952797
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.05168867829117
$ws.Range("D7").Value = 1.059043098018878
$ws.Range("E7").Value = 1.058507312948584
$ws.Range("F7").Value = 1.06905396146892
$ws.Range("I7").Value = 1.049590009823232
$ws.Range("J7").Value = 1.056144382627869
$ws.Range("K7").Value = 1.061469478095143
$ws.Range("L7").Value = 1.060934985644006
$ws.Range("M7").Value = 1.071456445437781
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050499776960109
$ws.Range("D8").Value = 1.058080491062007
$ws.Range("E8").Value = 1.057457547193261
$ws.Range("F8").Value = 1.06791732529144
$ws.Range("I8").Value = 1.049248624170786
$ws.Range("J8").Value = 1.055416008116475
$ws.Range("K8").Value = 1.060751700222901
$ws.Range("L8").Value = 1.060130420454089
$ws.Range("M8").Value = 1.070562534876355
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04840907777859
$ws.Range("D9").Value = 1.056387675566138
$ws.Range("E9").Value = 1.055614289921757
$ws.Range("F9").Value = 1.06592143399018
$ws.Range("I9").Value = 1.048638957386729
$ws.Range("J9").Value = 1.054131739797113
$ws.Range("K9").Value = 1.059485149646745
$ws.Range("L9").Value = 1.05871418348251
$ws.Range("M9").Value = 1.06898939721439
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047018247647768
$ws.Range("D10").Value = 1.055261527213209
$ws.Range("E10").Value = 1.054389962166734
$ws.Range("F10").Value = 1.064595652126908
$ws.Range("I10").Value = 1.048227109447058
$ws.Range("J10").Value = 1.053275105981339
$ws.Range("K10").Value = 1.058639689324034
$ws.Range("L10").Value = 1.057771117143967
$ws.Range("M10").Value = 1.067942099414913
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046416723841862
$ws.Range("D11").Value = 1.05477448068423
$ws.Range("E11").Value = 1.053860904946502
$ws.Range("F11").Value = 1.064022736639709
$ws.Range("I11").Value = 1.048047502237815
$ws.Range("J11").Value = 1.052904080627927
$ws.Range("K11").Value = 1.058273351594534
$ws.Range("L11").Value = 1.057363033994001
$ws.Range("M11").Value = 1.067488972044367
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046193399827564
$ws.Range("D12").Value = 1.054593659131273
$ws.Range("E12").Value = 1.05366455415484
$ws.Range("F12").Value = 1.063810105824603
$ws.Range("I12").Value = 1.047980597376858
$ws.Range("J12").Value = 1.052766251825586
$ws.Range("K12").Value = 1.058137241425455
$ws.Range("L12").Value = 1.057211495708183
$ws.Range("M12").Value = 1.067320715784948
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04624129864618
$ws.Range("D13").Value = 1.054632441905352
$ws.Range("E13").Value = 1.053706664590612
$ws.Range("F13").Value = 1.063855707886792
$ws.Range("I13").Value = 1.047994957316939
$ws.Range("J13").Value = 1.052795817162525
$ws.Range("K13").Value = 1.058166439129561
$ws.Range("L13").Value = 1.057243999265268
$ws.Range("M13").Value = 1.067356804793129
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046398261583891
$ws.Range("D14").Value = 1.054759532089052
$ws.Range("E14").Value = 1.053844671161203
$ws.Range("F14").Value = 1.064005156920053
$ws.Range("I14").Value = 1.048041975748205
$ws.Range("J14").Value = 1.052892687924205
$ws.Range("K14").Value = 1.058262101413407
$ws.Range("L14").Value = 1.057350506928515
$ws.Range("M14").Value = 1.067475062790723
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046494986074845
$ws.Range("D15").Value = 1.054837848415366
$ws.Range("E15").Value = 1.053929723426451
$ws.Range("F15").Value = 1.06409726070923
$ws.Range("I15").Value = 1.048070920100971
$ws.Range("J15").Value = 1.052952371447509
$ws.Range("K15").Value = 1.058321037366224
$ws.Range("L15").Value = 1.057416135414284
$ws.Range("M15").Value = 1.067547932832196
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04705818397439
$ws.Range("D16").Value = 1.055293863310257
$ws.Range("E16").Value = 1.054425096947098
$ws.Range("F16").Value = 1.064633699159399
$ws.Range("I16").Value = 1.048239002597633
$ws.Range("J16").Value = 1.053299727760758
$ws.Range("K16").Value = 1.058663996834276
$ws.Range("L16").Value = 1.057798206093431
$ws.Range("M16").Value = 1.067972179680709
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047411655442684
$ws.Range("D17").Value = 1.055580066661255
$ws.Range("E17").Value = 1.054736123253007
$ws.Range("F17").Value = 1.064970503743199
$ws.Range("I17").Value = 1.048344095630778
$ws.Range("J17").Value = 1.053517589968735
$ws.Range("K17").Value = 1.058879060749371
$ws.Range("L17").Value = 1.058037942341203
$ws.Range("M17").Value = 1.068238395762395
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047617898251002
$ws.Range("D18").Value = 1.055747060390702
$ws.Range("E18").Value = 1.054917644190344
$ws.Range("F18").Value = 1.065167067659303
$ws.Range("I18").Value = 1.048405271616679
$ws.Range("J18").Value = 1.053644655847801
$ws.Range("K18").Value = 1.059004479899301
$ws.Range("L18").Value = 1.058177802480607
$ws.Range("M18").Value = 1.06839370968942
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047688233362288
$ws.Range("D19").Value = 1.055804010413342
$ws.Range("E19").Value = 1.054979555823516
$ws.Range("F19").Value = 1.065234109735496
$ws.Range("I19").Value = 1.048426110162118
$ws.Range("J19").Value = 1.053687980389311
$ws.Range("K19").Value = 1.059047240481894
$ws.Range("L19").Value = 1.058225495552042
$ws.Range("M19").Value = 1.068446673549065
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047373724165138
$ws.Range("D20").Value = 1.05554935393013
$ws.Range("E20").Value = 1.054702742247416
$ws.Range("F20").Value = 1.064934356253467
$ws.Range("I20").Value = 1.048332832860504
$ws.Range("J20").Value = 1.053494216385953
$ws.Range("K20").Value = 1.058855988899618
$ws.Range("L20").Value = 1.058012218207742
$ws.Range("M20").Value = 1.068209829708202
$ws.Range("B21").Value = 1.019999999999999
$ws.Range("C21").Value = 1.046352036899408
$ws.Range("D21").Value = 1.054722104721984
$ws.Range("E21").Value = 1.0538040271085
$ws.Range("F21").Value = 1.063961143046133
$ws.Range("I21").Value = 1.04802813525194
$ws.Range("J21").Value = 1.052864162261856
$ws.Range("K21").Value = 1.058233932239894
$ws.Range("L21").Value = 1.057319141907561
$ws.Range("M21").Value = 1.067440237216683
$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.045710291182707
$ws.Range("D22").Value = 1.054202497267394
$ws.Range("E22").Value = 1.053239922461246
$ws.Range("F22").Value = 1.063350261544883
$ws.Range("I22").Value = 1.047835456524006
$ws.Range("J22").Value = 1.052467944884127
$ws.Range("K22").Value = 1.057842612697048
$ws.Range("L22").Value = 1.056883620766571
$ws.Range("M22").Value = 1.066956685223589
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046050432608052
$ws.Range("D23").Value = 1.054477901473633
$ws.Range("E23").Value = 1.05353887410498
$ws.Range("F23").Value = 1.063674004485854
$ws.Range("I23").Value = 1.047937703501518
$ws.Range("J23").Value = 1.052677994161083
$ws.Range("K23").Value = 1.058050077886957
$ws.Range("L23").Value = 1.057114475260624
$ws.Range("M23").Value = 1.067212994352018
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047390863464025
$ws.Range("D24").Value = 1.055563231517983
$ws.Range("E24").Value = 1.054717825364109
$ws.Range("F24").Value = 1.064950689403807
$ws.Range("I24").Value = 1.048337922401718
$ws.Range("J24").Value = 1.053504777931264
$ws.Range("K24").Value = 1.058866414149672
$ws.Range("L24").Value = 1.058023841756049
$ws.Range("M24").Value = 1.06822273737176
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048949056165239
$ws.Range("D25").Value = 1.056824893730156
$ws.Range("E25").Value = 1.056090028222428
$ws.Range("F25").Value = 1.066436578943907
$ws.Range("I25").Value = 1.048797526458343
$ws.Range("J25").Value = 1.054463839187443
$ws.Range("K25").Value = 1.059812781355032
$ws.Range("L25").Value = 1.059080127888
$ws.Range("M25").Value = 1.069395839739697
